# Rename the ResourceCreator.* header labels on Sheet1 row 1 so that the
# "institution" columns (E/F) are qualified as ResourceCreatorInstitution.*
# and the "person" columns (G/H/I/J/K) are qualified as ResourceCreatorPerson.*
# Also move the active selection from A4 to E2 (matches the merged-in sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "ResourceCreatorInstitution.Institution.name"
$ws.Range("F1").Value = "ResourceCreatorInstitution.role"
$ws.Range("G1").Value = "ResourceCreatorPerson.Person.lastName"
$ws.Range("H1").Value = "ResourceCreatorPerson.Person.firstName"
$ws.Range("I1").Value = "ResourceCreatorPerson.Person.email"
$ws.Range("J1").Value = "ResourceCreatorPerson.Person.Institution.name"
$ws.Range("K1").Value = "ResourceCreatorPerson.role"

$ws.Activate()
$ws.Range("E2").Select() | Out-Null
